# Repull data, push all data, mean calculation
# Updates column F (dSF) values for rows 2-31 to reflect newly pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 7
    3  = -3
    4  = -2
    5  = -1
    6  = -5
    7  = 8
    8  = -6
    9  = 2
    10 = 6
    12 = -7
    15 = 8
    16 = 10
    17 = -3
    19 = -3
    20 = -8
    21 = -9
    22 = -1
    23 = 4
    24 = -8
    25 = -2
    26 = 3
    28 = -8
    29 = -3
    31 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
